$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format first so numeric-looking strings (e.g. "1.00", "8.30")
# are preserved exactly instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "76.513.99"
$ws.Range("E2").Value = "  -0.63%  "

# Row 3
$ws.Range("D3").Value = "3.076.55"
$ws.Range("E3").Value = "  +3.69%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "198.51"
$ws.Range("E5").Value = "  -1.49%  "

# Row 6
$ws.Range("D6").Value = "616.72"
$ws.Range("E6").Value = "  +3.06%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.22%  "

# Row 9
$ws.Range("D9").Value = "0.208"
$ws.Range("E9").Value = "  +6.12%  "

# Row 10
$ws.Range("D10").Value = "3.073.67"
$ws.Range("E10").Value = "  +3.70%  "

# Row 11
$ws.Range("D11").Value = "0.439"
$ws.Range("E11").Value = "  -1.78%  "

# Row 12
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
$ws.Range("D13").Value = "5.23"
$ws.Range("E13").Value = "  +6.45%  "

# Row 14
$ws.Range("D14").Value = "3.609.67"
$ws.Range("E14").Value = "  +2.88%  "

# Row 15
$ws.Range("D15").Value = "28.96"
$ws.Range("E15").Value = "  +2.05%  "

# Row 16
$ws.Range("D16").Value = "76.449.70"
$ws.Range("E16").Value = "  -0.55%  "

# Row 17
$ws.Range("D17").Value = "0.0000193"
$ws.Range("E17").Value = "  +2.52%  "

# Row 18
$ws.Range("D18").Value = "3.064.55"
$ws.Range("E18").Value = "  +3.80%  "

# Row 19
$ws.Range("D19").Value = "13.58"
$ws.Range("E19").Value = "  +0.36%  "

# Row 20
$ws.Range("D20").Value = "8.92"
$ws.Range("E20").Value = "  +2.22%  "

# Row 21
$ws.Range("D21").Value = "381.25"
$ws.Range("E21").Value = "  +1.94%  "

# Row 22
$ws.Range("D22").Value = "2.46"
$ws.Range("E22").Value = "  +8.74%  "

# Row 23
$ws.Range("D23").Value = "4.40"
$ws.Range("E23").Value = "  +1.41%  "

# Row 24
$ws.Range("D24").Value = "3.224.50"
$ws.Range("E24").Value = "  +3.52%  "

# Row 25
$ws.Range("D25").Value = "72.49"
$ws.Range("E25").Value = "  -0.40%  "

# Row 26
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$ws.Range("D27").Value = "4.33"
$ws.Range("E27").Value = "  +1.15%  "

# Row 28
$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  +1.41%  "

# Row 29
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  -0.03%  "

# Row 30
$ws.Range("E30").Value = "  +0.20%  "

# Row 31
$ws.Range("D31").Value = "8.30"
$ws.Range("E31").Value = "  +0.03%  "

# Row 32
$ws.Range("E32").Value = "  +1.20%  "

# Row 33
$ws.Range("D33").Value = "498.41"
$ws.Range("E33").Value = "  -0.18%  "

# Row 34
$ws.Range("E34").Value = "  +4.36%  "

# Row 35
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.00%  "

# Row 36
$ws.Range("E36").Value = "  +12.92%  "

# Row 37
$ws.Range("D37").Value = "20.67"
$ws.Range("E37").Value = "  +2.24%  "

# Row 38
$ws.Range("D38").Value = "162.74"
$ws.Range("E38").Value = "  -2.14%  "

# Row 39
$ws.Range("D39").Value = "20.05"
$ws.Range("E39").Value = "  +1.28%  "

# Row 40
$ws.Range("D40").Value = "193.02"
$ws.Range("E40").Value = "  +6.43%  "

# Row 41
$ws.Range("D41").Value = "0.379"
$ws.Range("E41").Value = "  -5.13%  "

# Row 42
$ws.Range("E42").Value = "  -8.32%  "

# Row 43
$ws.Range("E43").Value = "  +0.03%  "

# Row 44
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.790"
$ws.Range("E44").Value = "  +19.29%  "

# Row 45
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D45").Value = "5.11"
$ws.Range("E45").Value = "  +3.53%  "

# Row 46
$ws.Range("D46").Value = "1.24"
$ws.Range("E46").Value = "  +4.35%  "

# Row 47
$ws.Range("D47").Value = "41.29"
$ws.Range("E47").Value = "  +2.91%  "

# Row 48
$ws.Range("E48").Value = "  -0.72%  "

# Row 49
$ws.Range("D49").Value = "2.43"
$ws.Range("E49").Value = "  +4.37%  "

# Row 50
$ws.Range("D50").Value = "0.596"
$ws.Range("E50").Value = "  +0.42%  "

# Row 51
$ws.Range("D51").Value = "3.87"
$ws.Range("E51").Value = "  -0.82%  "
